$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-08 Wednesday", "2025-10-09 Thursday"),
    @("849÷5=", "430÷4="),
    @("348÷9=", "389÷7="),
    @("945÷5=", "656÷6="),
    @("229÷9=", "746÷6="),
    @("810÷6=", "951÷2="),
    @("566÷6=", "443÷8="),
    @("445÷8=", "459÷8="),
    @("125÷7=", "479÷3="),
    @("563÷4=", "713÷9="),
    @("151÷4=", "603÷7="),
    @("199÷9=", "822÷3="),
    @("226÷9=", "153÷9="),
    @("113÷8=", "272÷7="),
    @("977÷2=", "530÷9="),
    @("392÷9=", "914÷3="),
    @("220÷5=", "877÷7="),
    @("963÷8=", "648÷8="),
    @("741÷8=", "866÷8="),
    @("639÷9=", "219÷2="),
    @("494÷4=", "426÷9="),
    @("776÷7=", "895÷3="),
    @("562÷4=", "738÷5="),
    @("119÷9=", "695÷7="),
    @("768÷3=", "252÷9="),
    @("105÷5=", "726÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
